$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-11 (columns A item id, B gene(s) - may be blank, E always "no")
$data = @(
    @("cg08369368", "NSD1"),
    @("cg03890691", "DOC2A"),
    @("cg15772157", ""),
    @("cg24789467", "SHROOM1"),
    @("cg02716826", "SUGT1P1;AQP3"),
    @("cg27151362", "DOC2A"),
    @("cg26612727", "ZPBP2"),
    @("cg23960707", "SFRS1"),
    @("cg26276120", "TPI1"),
    @("cg18310639", "CAPN2")
)

$ws.Range("E1").Value = "reverses"

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 5).Value = "no"
}

$ws.Range("C21").Select()
